$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data: insert a header row above, and add a new first column
# Final layout:
#   A1 = "What is it?"   B1 = "How much?"   C1 = "Units?"     (bold headers)
#   A2 = "2m symbol rate" B2 = 9600          C2 = "bps"        (left aligned)

$ws.Range("A2").Value = "2m symbol rate"
$ws.Range("B2").Value = 9600

$ws.Range("B1").Value = "How much?"
$ws.Range("C1").Value = "Units?"
$ws.Range("A1").Value = "What is it?"

$ws.Range("C2").Value = "bps"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 15.6640625
$ws.Columns.Item(2).ColumnWidth = 13
$ws.Columns.Item(3).ColumnWidth = 14.83203125

# Alignment: whole used area left aligned
$ws.Range("A1:C2").HorizontalAlignment = -4131

# Bold the header row
$ws.Range("A1:C1").Font.Bold = $true

$ws.Range("B6").Select()
